$wb = $excel.ActiveWorkbook

# Sheets: 1=geometry, 2=kpt.alex, 3=List1
$wsAlex = $wb.Worksheets.Item("kpt.alex")
$wsList = $wb.Worksheets.Item("List1")

# --- sheet3 (List1): add two new rows with new values/strings ---
$wsList.Range("A12").Value = "MM_LOAD.LZ"
$wsList.Range("B12").Value = 402
$wsList.Range("C12").Value = "load position screen"

$wsList.Range("A13").Value = "OPT_BAR.LZ"
$wsList.Range("B13").Value = 10
$wsList.Range("C13").Value = "looks like chunk of window frame?"

# --- selections / active sheet ---
# kpt.alex (sheet2): change selection to F19, no longer tabSelected
$wsAlex.Range("F19").Select()

# List1 (sheet3): change selection to C14, becomes tabSelected
$wsList.Range("C14").Select()
$wsList.Activate()

$wb.Save()
